$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "G2" = 0.208324
    "H2" = 0.624972
    "I2" = 0.002558447967581038
    "J2" = 0.002558447967581039
    "M2" = 3.997860333333333
    "N2" = 11.993581
    "O2" = 0.6810627835305383
    "P2" = 0.6810627835305384
    "Q2" = 0.8328502560813332
    "R2" = 7.495652304731999
    "S2" = 0.00174246369431879
    "T2" = 0.001742463694318791
    "G3" = 0.208324
    "H3" = 0.624972
    "I3" = 0.002558447967581038
    "J3" = 0.002558447967581039
    "O3" = 0.07471986524093907
    "P3" = 0.07471986524093907
    "Q3" = 0.09137257299200001
    "R3" = 0.822353156928
    "S3" = 0.0001911668873636096
    "T3" = 0.0001911668873636097
    "G4" = 0.208324
    "H4" = 0.624972
    "I4" = 0.002558447967581038
    "J4" = 0.002558447967581039
    "M4" = 1.433563666666667
    "N4" = 4.300691
    "O4" = 0.2442173512285226
    "P4" = 0.2442173512285226
    "Q4" = 0.2986457172946667
    "R4" = 2.687811455652
    "S4" = 0.0006248173858986382
    "T4" = 0.0006248173858986383
    "I5" = 0.6383194838544957
    "J5" = 0.6383194838544958
    "M5" = 3.997860333333333
    "N5" = 11.993581
    "O5" = 0.6810627835305383
    "P5" = 0.6810627835305384
    "Q5" = 207.7918145400321
    "R5" = 1870.126330860289
    "S5" = 0.4347356444557193
    "T5" = 0.4347356444557194
    "I6" = 0.6383194838544957
    "J6" = 0.6383194838544958
    "O6" = 0.07471986524093907
    "P6" = 0.07471986524093907
    "S6" = 0.0476951458142737
    "T6" = 0.04769514581427371
    "I7" = 0.6383194838544957
    "J7" = 0.6383194838544958
    "M7" = 1.433563666666667
    "N7" = 4.300691
    "O7" = 0.2442173512285226
    "P7" = 0.2442173512285226
    "Q7" = 74.51055582698656
    "R7" = 670.595002442879
    "S7" = 0.1558886935845026
    "T7" = 0.1558886935845026
    "G8" = 25.552936
    "H8" = 76.65880799999999
    "I8" = 0.3138181734938286
    "J8" = 0.3138181734938286
    "M8" = 3.997860333333333
    "N8" = 11.993581
    "O8" = 0.6810627835305383
    "P8" = 0.6810627835305384
    "Q8" = 102.1570692346053
    "R8" = 919.4136231114478
    "S8" = 0.2137298787621763
    "T8" = 0.2137298787621763
    "G9" = 25.552936
    "H9" = 76.65880799999999
    "I9" = 0.3138181734938286
    "J9" = 0.3138181734938286
    "O9" = 0.07471986524093907
    "P9" = 0.07471986524093907
    "Q9" = 11.207722153088
    "R9" = 100.869499377792
    "S9" = 0.02344845163361651
    "T9" = 0.02344845163361651
    "G10" = 25.552936
    "H10" = 76.65880799999999
    "I10" = 0.3138181734938286
    "J10" = 0.3138181734938286
    "M10" = 1.433563666666667
    "N10" = 4.300691
    "O10" = 0.2442173512285226
    "P10" = 0.2442173512285226
    "Q10" = 36.63176062625867
    "R10" = 329.685845636328
    "S10" = 0.07663984309803579
    "T10" = 0.07663984309803579
    "G11" = 3.688911666666666
    "H11" = 11.066735
    "I11" = 0.04530389468409456
    "J11" = 0.04530389468409456
    "M11" = 3.997860333333333
    "N11" = 11.993581
    "O11" = 0.6810627835305383
    "P11" = 0.6810627835305384
    "Q11" = 14.74775362533722
    "R11" = 132.729782628035
    "S11" = 0.0308547966183238
    "T11" = 0.03085479661832381
    "G12" = 3.688911666666666
    "H12" = 11.066735
    "I12" = 0.04530389468409456
    "J12" = 0.04530389468409456
    "O12" = 0.07471986524093907
    "P12" = 0.07471986524093907
    "Q12" = 1.617986168293333
    "R12" = 14.56187551464
    "S12" = 0.003385100905685241
    "T12" = 0.003385100905685241
    "G13" = 3.688911666666666
    "H13" = 11.066735
    "I13" = 0.04530389468409456
    "J13" = 0.04530389468409456
    "M13" = 1.433563666666667
    "N13" = 4.300691
    "O13" = 0.2442173512285226
    "P13" = 0.2442173512285226
    "Q13" = 5.288289734876112
    "R13" = 47.594607613885
    "S13" = 0.01106399716008552
    "T13" = 0.01106399716008552
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
